$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.030001401901245
$ws.Range("B1").Value = 1.059765696525574
$ws.Range("C1").Value = 1.619904756546021
$ws.Range("D1").Value = 2.530822992324829
$ws.Range("E1").Value = 0.7748891115188599
